$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking identifiers that must stay text
# (e.g. "1.00", "6.70", "64.100.88") so force text format before assignment.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.100.88"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.476.66"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.54"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.51"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.481"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("E10").Value = "  -1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.387"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.068.80"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.478.30"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.138.65"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.33"
$ws.Range("E17").Value = "  -6.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.95"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.41"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "384.14"
$ws.Range("E21").Value = "  -2.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.571"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.620.05"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.70"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -2.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.12"
$ws.Range("E30").Value = "  -4.03%  "
$ws.Range("E31").Value = "  -4.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").Value = "  -4.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.509.87"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.151"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.75"
$ws.Range("E38").Value = "  -2.44%  "
$ws.Range("E39").Value = "  -4.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.71"
$ws.Range("E40").Value = "  +0.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0774"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.38"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.62"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.56"
$ws.Range("E47").Value = "  -6.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.70"
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.901"
$ws.Range("E50").Value = "  +0.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.355.12"
$ws.Range("E51").Value = "  -4.85%  "
